$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "Índice"
$ws.Range("B1").Value = "Distancia"
$ws.Range("C1").Value = "max"
$ws.Range("D1").Value = "min"
$ws.Range("E1").Value = "Tempo"

# Data rows
$data = @(
    @(0, 1057, 1057, 1057, 0.01068532466888428),
    @(1, 1040, 1040, 1040, 0.01092848777770996),
    @(2, 973, 973, 973, 0.01275227069854736),
    @(3, 1224, 1224, 1224, 0.01249235471089681),
    @(4, 883, 883, 883, 0.01258847713470459),
    @(5, 1040, 1040, 1040, 0.01271528402964274),
    @(6, 1053, 1053, 1053, 0.01273972988128662),
    @(7, 957, 957, 957, 0.01048529148101807),
    @(8, 886, 886, 886, 0.01282922426859538),
    @(9, 1049, 1049, 1049, 0.01225011348724365)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $rowData = $data[$i]
    $ws.Cells.Item($row, 1).Value = $rowData[0]
    $ws.Cells.Item($row, 2).Value = $rowData[1]
    $ws.Cells.Item($row, 3).Value = $rowData[2]
    $ws.Cells.Item($row, 4).Value = $rowData[3]
    $ws.Cells.Item($row, 5).Value = $rowData[4]
}
